# DTY: Shipping DTY eligible order from workhouse
#
# Re-point the envt data sheet's "test21" environment to "test17":
#   - A2/C2/D2/F2/G2 text values move from the test21.* URLs / labels to the
#     equivalent test17.* ones.
#   - C2's hyperlink "display" text was already the test17 URL, so once the
#     cell text itself becomes test17 that cached display text is redundant
#     and Excel drops it from the saved hyperlink.
#   - The saved selection moves from A2 to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash original cell formatting in scratch cells far outside the used range,
# since rebuilding hyperlinks resets the touched cells to the generic
# "Hyperlink" built-in style.
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("Z3").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$ws.Range("J2").Copy()
$ws.Range("Z5").PasteSpecial(-4122)

# Deleting hyperlinks via any range clears the whole sheet's hyperlink collection
# in this runtime, so rebuild all 5 after deleting once.
$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "https://test17.cliotest.com/cabicentral/control/main") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://sandbox.cabiclio.com/warehouse/control/main", $null, $null, "https://sandbox.cabiclio.com/warehouse/control/main") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://mirandakate.cabitest5.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "https://test18.cliotest.com/backoffice/control/main", $null, $null, "https://test18.cliotest.com/backoffice/control/main") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J2"), "http://webmail.cabiclio.com/") | Out-Null

# Re-assert the actual cell text (Hyperlinks.Add's TextToDisplay, when given,
# would otherwise overwrite the cell text with the display string).
$ws.Range("A2").Value = "https://test17.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test17.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test17.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest17"
$ws.Range("G2").Value = "test17"

# Restore original formatting from the scratch cells
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("Z4").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("Z5").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# Clean up scratch cells
$ws.Range("Z1:Z5").Clear()

# Update the active selection to C12
$ws.Range("C12").Select() | Out-Null
